{"js": "const replacements = [\n  [\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"],\n  [\"68\u00d730=2040\", \"35\u00d740=1400\"],\n  [\"12\u00d745=540\", \"43\u00d799=4257\"],\n  [\"88\u00d726=2288\", \"49\u00d795=4655\"],\n  [\"53\u00d767=3551\", \"82\u00d774=6068\"],\n  [\"65\u00d793=6045\", \"94\u00d726=2444\"],\n  [\"63\u00d732=2016\", \"54\u00d735=1890\"],\n  [\"62\u00d718=1116\", \"43\u00d776=3268\"],\n  [\"22\u00d786=1892\", \"87\u00d757=4959\"],\n  [\"41\u00d775=3075\", \"13\u00d791=1183\"],\n  [\"76\u00d772=5472\", \"92\u00d790=8280\"],\n  [\"53\u00d725=1325\", \"99\u00d789=8811\"],\n  [\"71\u00d785=6035\", \"80\u00d742=3360\"],\n  [\"67\u00d758=3886\", \"38\u00d741=1558\"],\n  [\"33\u00d755=1815\", \"15\u00d726=390\"],\n  [\"92\u00d781=7452\", \"26\u00d786=2236\"],\n  [\"36\u00d763=2268\", \"44\u00d726=1144\"],\n  [\"96\u00d767=6432\", \"55\u00d718=990\"],\n  [\"40\u00d738=1520\", \"94\u00d749=4606\"],\n  [\"12\u00d722=264\", \"61\u00d777=4697\"],\n  [\"74\u00d741=3034\", \"60\u00d794=5640\"],\n  [\"89\u00d719=1691\", \"67\u00d794=6298\"],\n  [\"47\u00d755=2585\", \"80\u00d741=3280\"],\n  [\"44\u00d725=1100\", \"68\u00d773=4964\"],\n  [\"90\u00d794=8460\", \"31\u00d777=2387\"],\n  [\"21\u00d747=987\", \"96\u00d752=4992\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"),\n    @(\"68\u00d730=2040\", \"35\u00d740=1400\"),\n    @(\"12\u00d745=540\", \"43\u00d799=4257\"),\n    @(\"88\u00d726=2288\", \"49\u00d795=4655\"),\n    @(\"53\u00d767=3551\", \"82\u00d774=6068\"),\n    @(\"65\u00d793=6045\", \"94\u00d726=2444\"),\n    @(\"63\u00d732=2016\", \"54\u00d735=1890\"),\n    @(\"62\u00d718=1116\", \"43\u00d776=3268\"),\n    @(\"22\u00d786=1892\", \"87\u00d757=4959\"),\n    @(\"41\u00d775=3075\", \"13\u00d791=1183\"),\n    @(\"76\u00d772=5472\", \"92\u00d790=8280\"),\n    @(\"53\u00d725=1325\", \"99\u00d789=8811\"),\n    @(\"71\u00d785=6035\", \"80\u00d742=3360\"),\n    @(\"67\u00d758=3886\", \"38\u00d741=1558\"),\n    @(\"33\u00d755=1815\", \"15\u00d726=390\"),\n    @(\"92\u00d781=7452\", \"26\u00d786=2236\"),\n    @(\"36\u00d763=2268\", \"44\u00d726=1144\"),\n    @(\"96\u00d767=6432\", \"55\u00d718=990\"),\n    @(\"40\u00d738=1520\", \"94\u00d749=4606\"),\n    @(\"12\u00d722=264\", \"61\u00d777=4697\"),\n    @(\"74\u00d741=3034\", \"60\u00d794=5640\"),\n    @(\"89\u00d719=1691\", \"67\u00d794=6298\"),\n    @(\"47\u00d755=2585\", \"80\u00d741=3280\"),\n    @(\"44\u00d725=1100\", \"68\u00d773=4964\"),\n    @(\"90\u00d794=8460\", \"31\u00d777=2387\"),\n    @(\"21\u00d747=987\", \"96\u00d752=4992\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
